# Add "attackSpeed|Float" column to the ActorTable sheet (Actor.xlsx),
# mirroring the same column already present on ActorPowerLevelTable.

$wb = $excel.ActiveWorkbook

$wsActor = $wb.Worksheets.Item("ActorTable")

# --- Header (row 1, column E) ---
$wsActor.Cells.Item(1, 5).Value = "attackSpeed|Float"

# --- Data rows (2..18) all get value 1, same as column D ---
for ($r = 2; $r -le 18; $r++) {
    $wsActor.Cells.Item($r, 5).Value = 1
}

# Column E width, matching the sibling sheet's look-and-feel
$wsActor.Columns.Item(5).ColumnWidth = 15.57

# ActorTable becomes the active/selected sheet (tabSelected moves from
# ActorPowerLevelTable to ActorTable).
$wsActor.Activate()
